$d = $word.ActiveDocument

# --- Locate the "Static Fields" heading paragraph (the bold run that
# directly follows the inline picture) and split it into two paragraphs:
#   1) an empty paragraph whose mark now carries the bold run formatting
#   2) the original paragraph mark formatting (non-bold) + the bold
#      "Static Fields" run, now carrying the <w:lastRenderedPageBreak/>
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Static Fields" + [char]13) {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the 'Static Fields' heading paragraph"
}

$xmlns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$newXml = '<w:p ' + $xmlns + '>' +
    '<w:pPr>' +
        '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
        '<w:spacing w:beforeAutospacing="1" w:after="0" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/>' +
            '<w:b/>' +
            '<w:bCs/>' +
            '<w:color w:val="2D2F31"/>' +
            '<w:kern w:val="0"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w14:ligatures w14:val="none"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
'</w:p>' +
'<w:p ' + $xmlns + '>' +
    '<w:pPr>' +
        '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
        '<w:spacing w:beforeAutospacing="1" w:after="0" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/>' +
            '<w:color w:val="2D2F31"/>' +
            '<w:kern w:val="0"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w14:ligatures w14:val="none"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/>' +
            '<w:b/>' +
            '<w:bCs/>' +
            '<w:color w:val="2D2F31"/>' +
            '<w:kern w:val="0"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w14:ligatures w14:val="none"/>' +
        '</w:rPr>' +
        '<w:lastRenderedPageBreak/>' +
        '<w:t>Static Fields</w:t>' +
    '</w:r>' +
'</w:p>'

$targetPara.Range.InsertXML($newXml)

# --- Remove the now-redundant <w:lastRenderedPageBreak/> from the
# "Static fields are store outside the object." run. Use a range that
# stops one character short of the paragraph mark so only the run
# content is rewritten and the <w:p>/<w:pPr> (with their original
# rsid/paraId attributes) are left completely untouched.
$bodyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Static fields are store outside the object." + [char]13) {
        $bodyPara = $p
        break
    }
}
if ($bodyPara -eq $null) {
    throw "Could not find the 'Static fields are store outside the object.' paragraph"
}

$bodyRange = $d.Range($bodyPara.Range.Start, $bodyPara.Range.End - 1)

$bodyRunXml = '<w:p ' + $xmlns + '>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/>' +
            '<w:color w:val="2D2F31"/>' +
            '<w:kern w:val="0"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w14:ligatures w14:val="none"/>' +
        '</w:rPr>' +
        '<w:t>Static fields are store outside the object.</w:t>' +
    '</w:r>' +
'</w:p>'

$bodyRange.InsertXML($bodyRunXml)
